$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "Listing"
$ws.Range("G6").Value = "Listing"
$ws.Range("G9").Value = "Listing, Property"
$ws.Range("G13").Value = "Listing"
$ws.Range("G14").Value = "Listing"
$ws.Range("G16").Value = "Listing, Property, host"
$ws.Range("G22").Value = "Listing, Property, host"
$ws.Range("G30").Value = "Listing"

$ws.Rows.Item(17).RowHeight = 16.8

$ws.Rows.Item(18).Select()
